$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: user_name, change "Yasuhiro Kato" -> "Tsubasa Endo" for rows 2-16 ---
foreach ($r in 2..16) {
    $ws.Range("C$r").Value = "Tsubasa Endo"
}

# --- Row 5: type operation -> error ---
$ws.Range("B5").Value = "error"

# --- Row 7: type error -> operation ---
$ws.Range("B7").Value = "operation"

# --- Column J: capimg filenames ---
$ws.Range("J2").Value = "bdot20240415_141954/1.png"
$ws.Range("J3").Value = "bdot20240415_141954/2.png"
$ws.Range("J4").Value = "bdot20240415_141954/3.png"
$ws.Range("J5").Value = "bdot20240415_141954/4.png"
$ws.Range("J6").Value = "bdot20240415_141954/5.png"
$ws.Range("J7").Value = "bdot20240415_141954/5.png"
$ws.Range("J8").Value = "bdot20240415_141954/6.png"
$ws.Range("J9").Value = "bdot20240415_141954/7.png"
$ws.Range("J10").Value = "bdot20240415_141954/8.png"
$ws.Range("J11").Value = "bdot20240415_141954/9.png"
$ws.Range("J12").Value = "bdot20240415_141954/10.png"
$ws.Range("J13").Value = "bdot20240415_141954/1.png"
$ws.Range("J14").Value = "bdot20240415_141954/2.png"
$ws.Range("J15").Value = "bdot20240415_141954/3.png"
$ws.Range("J16").Value = "bdot20240415_141954/11.png"

# --- Column K: explanation text ---
$ws.Range("K2").Value = "「スタート」ボタンをクリックする"
$ws.Range("K3").Value = "メニューから「設定」アイコンをクリックする"
$ws.Range("K4").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
$ws.Range("K5").Value = "0x80240fff エラー"
$ws.Range("K6").Value = "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"
$ws.Range("K7").Value = "メニューからターミナル(管理者)をクリックする"
$ws.Range("K8").Value = "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"
$ws.Range("K9").Value = "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"
$ws.Range("K10").Value = "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"
$ws.Range("K11").Value = "netsh winhttp show proxy と入力し、[Enter]キーを押す"
$ws.Range("K12").Value = "netsh winhttp reset proxy と入力し、[Enter]キーを押す"
$ws.Range("K13").Value = "「スタート」ボタンをクリックする"
$ws.Range("K14").Value = "メニューから「設定」アイコンをクリックする"
$ws.Range("K15").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
$ws.Range("K16").Value = "「更新プログラムのチェック」ボタンをクリックする"

# --- Row 5 gains error_type / error_content (previously on row 7) ---
$ws.Range("L5").Value = "Error W"
$ws.Range("M5").Value = " エラーの Windows"

# --- Row 7 loses error_type / error_content (moved to row 5) ---
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
